$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 5
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = 3
$ws.Range("F9").Value = 0
$ws.Range("F13").Value = 2
$ws.Range("F15").Value = 2
$ws.Range("F18").Value = 3
$ws.Range("F21").Value = -2
$ws.Range("F27").Value = -2
$ws.Range("F32").Value = -2
$ws.Range("F34").Value = -2
$ws.Range("F35").Value = -1
$ws.Range("F42").Value = -2
$ws.Range("F50").Value = 1
$ws.Range("F51").Value = -5
$ws.Range("F53").Value = -1
$ws.Range("F57").Value = 0
